$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new columns C..G, replacing old Name/salary/remarks columns ---
# Order matters for shared-string index assignment, so set values in the same
# order they first appear in the target workbook.
$ws.Range("C1").Value = "Project reviews"
$ws.Range("C2").Value = "Sprint1"
$ws.Range("D1").Value = "review comments"
$ws.Range("E1").Value = "action plan"
$ws.Range("F1").Value = "owner"
$ws.Range("G1").Value = "status remarks"

# --- Date column B, rows 2..33: 2024-01-01 (45292) through 2024-02-01 (45323) ---
for ($i = 2; $i -le 33; $i++) {
    $serial = 45292 + ($i - 2)
    $cell = $ws.Cells.Item($i, 2)
    $cell.Value = $serial
    $cell.NumberFormat = "d-mmm-yy"
}

# --- Autofit the populated columns to better match column widths ---
$ws.Columns("B:E").AutoFit()
$ws.Columns("G:G").AutoFit()

# --- Final selection: entire column H (matches saved view state) ---
$ws.Columns("H:H").Select()
